$wb = $excel.ActiveWorkbook

# --- Sheet: Homeostatic ---
$ws = $wb.Worksheets.Item("Homeostatic")
$ws.Range("A2:D33").ClearContents()
$ws.Range("A2").Value = "Siglech"
$ws.Range("B2").Value = "Imm.M0Like.0"
$ws.Range("C2").Value = 0.5900871753692627
$ws.Range("D2").Value = 0.3601895734597156
$ws.Range("A3").Value = "Csf3r"
$ws.Range("B3").Value = "Imm.M0Like.0"
$ws.Range("C3").Value = 0.3081689476966858
$ws.Range("D3").Value = 0.6966824644549763
$ws.Range("A4").Value = "Csf1r"
$ws.Range("B4").Value = "Imm.M0Like.0"
$ws.Range("C4").Value = 0.2840373516082764
$ws.Range("D4").Value = 0.6042654028436019
$ws.Range("A5").Value = "Hexb"
$ws.Range("B5").Value = "Imm.M0Like.0"
$ws.Range("C5").Value = 0.2203921973705292
$ws.Range("D5").Value = 0.8222748815165877
$ws.Range("A6").Value = "Ly86"
$ws.Range("B6").Value = "Imm.M0Like.0"
$ws.Range("C6").Value = 0.244189664721489
$ws.Range("D6").Value = 0.7274881516587678
$ws.Range("A7").Value = "Cx3cr1"
$ws.Range("B7").Value = "Imm.M0Like.0"
$ws.Range("C7").Value = 0.1870851963758469
$ws.Range("D7").Value = 0.5734597156398105
$ws.Range("A8").Value = "P2ry12"
$ws.Range("B8").Value = "Imm.M0Like.0"
$ws.Range("C8").Value = 0.3282490670681
$ws.Range("D8").Value = 0.4360189573459716
$ws.Range("A9").Value = "Mertk"
$ws.Range("B9").Value = "Imm.M0Like.0"
$ws.Range("C9").Value = 0.04550057277083397
$ws.Range("D9").Value = 0.7914691943127962
$ws.Range("A10").Value = "Ctss"
$ws.Range("B10").Value = "Imm.M0Like.0"
$ws.Range("C10").Value = -0.1872172802686691
$ws.Range("D10").Value = 0.495260663507109
$ws.Range("A11").Value = "Cx3cr1"
$ws.Range("B11").Value = "Imm.M0Like.1"
$ws.Range("C11").Value = 0.7888805866241455
$ws.Range("D11").Value = 0.6593406593406593
$ws.Range("A12").Value = "Mertk"
$ws.Range("B12").Value = "Imm.M0Like.1"
$ws.Range("C12").Value = 0.4895345866680145
$ws.Range("D12").Value = 0.8434065934065934
$ws.Range("A13").Value = "Hexb"
$ws.Range("B13").Value = "Imm.M0Like.1"
$ws.Range("C13").Value = 0.454164057970047
$ws.Range("D13").Value = 0.8543956043956044
$ws.Range("A14").Value = "Csf3r"
$ws.Range("B14").Value = "Imm.M0Like.1"
$ws.Range("C14").Value = 0.4002310931682587
$ws.Range("D14").Value = 0.7335164835164835
$ws.Range("A15").Value = "Ly86"
$ws.Range("B15").Value = "Imm.M0Like.1"
$ws.Range("C15").Value = 0.3987973630428314
$ws.Range("D15").Value = 0.75
$ws.Range("A16").Value = "Csf1r"
$ws.Range("B16").Value = "Imm.M0Like.1"
$ws.Range("C16").Value = 0.3818937242031097
$ws.Range("D16").Value = 0.6373626373626373
$ws.Range("A17").Value = "P2ry12"
$ws.Range("B17").Value = "Imm.M0Like.1"
$ws.Range("C17").Value = 0.4483467638492584
$ws.Range("D17").Value = 0.4835164835164835
$ws.Range("A18").Value = "Siglech"
$ws.Range("B18").Value = "Imm.M0Like.1"
$ws.Range("C18").Value = 0.5787144303321838
$ws.Range("D18").Value = 0.3571428571428572
$ws.Range("A19").Value = "Ctss"
$ws.Range("B19").Value = "Imm.M0Like.1"
$ws.Range("C19").Value = 0.229726642370224
$ws.Range("D19").Value = 0.6016483516483516
$ws.Range("A20").Value = "Hexb"
$ws.Range("B20").Value = "Imm.M0Like.2"
$ws.Range("C20").Value = 0.3351731598377228
$ws.Range("D20").Value = 0.852017937219731
$ws.Range("A21").Value = "Siglech"
$ws.Range("B21").Value = "Imm.M0Like.2"
$ws.Range("C21").Value = 0.5694880485534668
$ws.Range("D21").Value = 0.3991031390134529
$ws.Range("A22").Value = "P2ry12"
$ws.Range("B22").Value = "Imm.M0Like.2"
$ws.Range("C22").Value = 0.4630667865276337
$ws.Range("D22").Value = 0.5112107623318386
$ws.Range("A23").Value = "Mertk"
$ws.Range("B23").Value = "Imm.M0Like.2"
$ws.Range("C23").Value = 0.2984011769294739
$ws.Range("D23").Value = 0.8385650224215246
$ws.Range("A24").Value = "Csf3r"
$ws.Range("B24").Value = "Imm.M0Like.2"
$ws.Range("C24").Value = 0.287078857421875
$ws.Range("D24").Value = 0.7085201793721974
$ws.Range("A25").Value = "Cx3cr1"
$ws.Range("B25").Value = "Imm.M0Like.2"
$ws.Range("C25").Value = 0.2182062566280365
$ws.Range("D25").Value = 0.5964125560538116
$ws.Range("A26").Value = "Ly86"
$ws.Range("B26").Value = "Imm.M0Like.2"
$ws.Range("C26").Value = 0.1176583170890808
$ws.Range("D26").Value = 0.7713004484304933
$ws.Range("A27").Value = "Ctss"
$ws.Range("B27").Value = "Imm.M0Like.2"
$ws.Range("C27").Value = 0.1129345670342445
$ws.Range("D27").Value = 0.5650224215246636
$ws.Range("A28").Value = "Csf1r"
$ws.Range("B28").Value = "Imm.M0Like.2"
$ws.Range("C28").Value = 0.07088322937488556
$ws.Range("D28").Value = 0.600896860986547
$ws.Rows("29:33").Delete()

# --- Sheet: MHCII ---
$ws = $wb.Worksheets.Item("MHCII")
$ws.Range("A2:D7").ClearContents()
$ws.Range("A2").Value = "Cd74"
$ws.Range("B2").Value = "Imm.MHCII.0"
$ws.Range("C2").Value = 5.328258991241455
$ws.Range("D2").Value = 0.6818181818181818
$ws.Range("A3").Value = "H2-Aa"
$ws.Range("B3").Value = "Imm.MHCII.0"
$ws.Range("C3").Value = 5.333882331848145
$ws.Range("D3").Value = 0.4318181818181818
$ws.Range("A4").Value = "H2-Ab1"
$ws.Range("B4").Value = "Imm.MHCII.0"
$ws.Range("C4").Value = 4.932478904724121
$ws.Range("D4").Value = 0.4318181818181818
$ws.Range("A5").Value = "H2-Eb1"
$ws.Range("B5").Value = "Imm.MHCII.0"
$ws.Range("C5").Value = 4.92715311050415
$ws.Range("D5").Value = 0.4318181818181818
$ws.Range("A6").Value = "Ciita"
$ws.Range("B6").Value = "Imm.MHCII.0"
$ws.Range("C6").Value = 6.172420024871826
$ws.Range("D6").Value = 0.4090909090909091
$ws.Rows("7:7").Delete()

# --- Sheet: Interferon ---
$ws = $wb.Worksheets.Item("Interferon")
$ws.Range("A2:D9").ClearContents()
$ws.Range("A2").Value = "Stat1"
$ws.Range("B2").Value = "Imm.Interferon.0"
$ws.Range("C2").Value = 2.863122940063477
$ws.Range("D2").Value = 0.717948717948718
$ws.Range("A3").Value = "Ifi209"
$ws.Range("B3").Value = "Imm.Interferon.0"
$ws.Range("C3").Value = 4.564236640930176
$ws.Range("D3").Value = 0.5641025641025641
$ws.Range("A4").Value = "Axl"
$ws.Range("B4").Value = "Imm.Interferon.0"
$ws.Range("C4").Value = 2.12287425994873
$ws.Range("D4").Value = 0.717948717948718
$ws.Range("A5").Value = "Ifi204"
$ws.Range("B5").Value = "Imm.Interferon.0"
$ws.Range("C5").Value = 3.268746376037598
$ws.Range("D5").Value = 0.5897435897435898
$ws.Range("A6").Value = "Stat2"
$ws.Range("B6").Value = "Imm.Interferon.0"
$ws.Range("C6").Value = 3.011718273162842
$ws.Range("D6").Value = 0.5384615384615384
$ws.Range("A7").Value = "Oasl2"
$ws.Range("B7").Value = "Imm.Interferon.0"
$ws.Range("C7").Value = 4.026572227478027
$ws.Range("D7").Value = 0.4615384615384616
$ws.Range("A8").Value = "Usp18"
$ws.Range("B8").Value = "Imm.Interferon.0"
$ws.Range("C8").Value = 3.895071506500244
$ws.Range("D8").Value = 0.4102564102564102
$ws.Rows("9:9").Delete()

# --- Sheet: DAM ---
$ws = $wb.Worksheets.Item("DAM")
$ws.Range("A2:D24").ClearContents()
$ws.Range("A2").Value = "Gpnmb"
$ws.Range("B2").Value = "Imm.DAM.0"
$ws.Range("C2").Value = 4.12397289276123
$ws.Range("D2").Value = 0.9230769230769231
$ws.Range("A3").Value = "Lyz2"
$ws.Range("B3").Value = "Imm.DAM.0"
$ws.Range("C3").Value = 2.719516038894653
$ws.Range("D3").Value = 0.6730769230769231
$ws.Range("A4").Value = "Lgals3"
$ws.Range("B4").Value = "Imm.DAM.0"
$ws.Range("C4").Value = 2.592077493667603
$ws.Range("D4").Value = 0.6346153846153846
$ws.Range("A5").Value = "Ctsb"
$ws.Range("B5").Value = "Imm.DAM.0"
$ws.Range("C5").Value = 1.726879477500916
$ws.Range("D5").Value = 0.8269230769230769
$ws.Range("A6").Value = "Apoe"
$ws.Range("B6").Value = "Imm.DAM.0"
$ws.Range("C6").Value = 1.632722139358521
$ws.Range("D6").Value = 0.8653846153846154
$ws.Range("A7").Value = "Ctsd"
$ws.Range("B7").Value = "Imm.DAM.0"
$ws.Range("C7").Value = 1.492807388305664
$ws.Range("D7").Value = 0.8269230769230769
$ws.Range("A8").Value = "Ftl1"
$ws.Range("B8").Value = "Imm.DAM.0"
$ws.Range("C8").Value = 2.148448467254639
$ws.Range("D8").Value = 0.3653846153846154
$ws.Range("A9").Value = "Spp1"
$ws.Range("B9").Value = "Imm.DAM.0"
$ws.Range("C9").Value = 1.38322651386261
$ws.Range("D9").Value = 0.4038461538461539
$ws.Range("A10").Value = "Fabp5"
$ws.Range("B10").Value = "Imm.DAM.0"
$ws.Range("C10").Value = 1.915103912353516
$ws.Range("D10").Value = 0.3076923076923077
$ws.Range("A11").Value = "Fth1"
$ws.Range("B11").Value = "Imm.DAM.0"
$ws.Range("C11").Value = 0.3970088064670563
$ws.Range("D11").Value = 0.4423076923076923
$ws.Range("A12").Value = "Lgals3"
$ws.Range("B12").Value = "Imm.DAM.1"
$ws.Range("C12").Value = 2.478747844696045
$ws.Range("D12").Value = 0.5447761194029851
$ws.Range("A13").Value = "Gpnmb"
$ws.Range("B13").Value = "Imm.DAM.1"
$ws.Range("C13").Value = 2.417728185653687
$ws.Range("D13").Value = 0.582089552238806
$ws.Range("A14").Value = "Apoe"
$ws.Range("B14").Value = "Imm.DAM.1"
$ws.Range("C14").Value = 1.087131977081299
$ws.Range("D14").Value = 0.7910447761194029
$ws.Range("A15").Value = "Ctsb"
$ws.Range("B15").Value = "Imm.DAM.1"
$ws.Range("C15").Value = 1.073116302490234
$ws.Range("D15").Value = 0.7014925373134329
$ws.Range("A16").Value = "Ctsd"
$ws.Range("B16").Value = "Imm.DAM.1"
$ws.Range("C16").Value = 0.9006564021110535
$ws.Range("D16").Value = 0.7761194029850746
$ws.Range("A17").Value = "Spp1"
$ws.Range("B17").Value = "Imm.DAM.1"
$ws.Range("C17").Value = 1.661806583404541
$ws.Range("D17").Value = 0.3880597014925373
$ws.Range("A18").Value = "Lpl"
$ws.Range("B18").Value = "Imm.DAM.1"
$ws.Range("C18").Value = 1.80559515953064
$ws.Range("D18").Value = 0.3358208955223881
$ws.Range("A19").Value = "Csf1"
$ws.Range("B19").Value = "Imm.DAM.1"
$ws.Range("C19").Value = 1.391711235046387
$ws.Range("D19").Value = 0.3059701492537313
$ws.Range("A20").Value = "Lyz2"
$ws.Range("B20").Value = "Imm.DAM.1"
$ws.Range("C20").Value = 1.102681040763855
$ws.Range("D20").Value = 0.3432835820895522
$ws.Range("A21").Value = "Fth1"
$ws.Range("B21").Value = "Imm.DAM.1"
$ws.Range("C21").Value = 0.3126291036605835
$ws.Range("D21").Value = 0.4776119402985075
$ws.Rows("22:24").Delete()

# --- Sheet: PVM ---
$ws = $wb.Worksheets.Item("PVM")
$ws.Range("A2:D4").ClearContents()
$ws.Range("A2").Value = "F13a1"
$ws.Range("B2").Value = "Imm.PVM.0"
$ws.Range("C2").Value = 6.067361831665039
$ws.Range("D2").Value = 0.7851239669421488
$ws.Range("A3").Value = "Mrc1"
$ws.Range("B3").Value = "Imm.PVM.0"
$ws.Range("C3").Value = 6.214034557342529
$ws.Range("D3").Value = 0.7107438016528925
$ws.Range("A4").Value = "Cd163"
$ws.Range("B4").Value = "Imm.PVM.0"
$ws.Range("C4").Value = 8.255138397216797
$ws.Range("D4").Value = 0.4049586776859504

# --- Sheet: Proliferative ---
$ws = $wb.Worksheets.Item("Proliferative")
$ws.Range("A2:D6").ClearContents()
$ws.Range("A2").Value = "Top2a"
$ws.Range("B2").Value = "Imm.Proliferative.0"
$ws.Range("C2").Value = 10.31621265411377
$ws.Range("D2").Value = 0.8333333333333334
$ws.Range("A3").Value = "Kif11"
$ws.Range("B3").Value = "Imm.Proliferative.0"
$ws.Range("C3").Value = 8.833640098571777
$ws.Range("D3").Value = 0.7777777777777778
$ws.Range("A4").Value = "Mki67"
$ws.Range("B4").Value = "Imm.Proliferative.0"
$ws.Range("C4").Value = 8.53589916229248
$ws.Range("D4").Value = 0.7222222222222222
$ws.Range("A5").Value = "Neil3"
$ws.Range("B5").Value = "Imm.Proliferative.0"
$ws.Range("C5").Value = 8.411818504333496
$ws.Range("D5").Value = 0.6666666666666666
$ws.Range("A6").Value = "Cenpf"
$ws.Range("B6").Value = "Imm.Proliferative.0"
$ws.Range("C6").Value = 8.655964851379395
$ws.Range("D6").Value = 0.6666666666666666
